# "Remove Facility from entire app"
# The "Facility" header (column C) and its sample data ("Facility name 1")
# are dropped from the Monthly Report sheet. Deleting the whole column
# shifts every later column one slot to the left and shrinks the sheet's
# used range from A1:AA2 down to A1:Z2, matching the data/shared-string
# changes in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Delete()

# Restore the active cell/selection on the sheet to reflect the
# post-edit state captured in the saved workbook.
[void]$ws.Range("E7").Select()
